$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 62, shifting existing rows 62-113 down to 67-118
$ws.Rows("62:66").Insert()

# Constant columns for this data block (Comercializadora del Agro de Limari / Durazno)
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria = "Durazno"

# Row 62
$ws.Cells.Item(62,1).Value2 = $mercadoId
$ws.Cells.Item(62,2).Value2 = $mercado
$ws.Cells.Item(62,3).Value2 = $region
$ws.Cells.Item(62,4).Value2 = 44629
$ws.Cells.Item(62,5).Value2 = $codreg
$ws.Cells.Item(62,6).Value2 = $tipo
$ws.Cells.Item(62,7).Value2 = $productoId
$ws.Cells.Item(62,8).Value2 = $producto
$ws.Cells.Item(62,9).Value2 = $categoriaId
$ws.Cells.Item(62,10).Value2 = $categoria
$ws.Cells.Item(62,11).Value2 = "Phillips Cling"
$ws.Cells.Item(62,12).Value2 = "Especial"
$ws.Cells.Item(62,13).Value2 = 16
$ws.Cells.Item(62,14).Value2 = 425000
$ws.Cells.Item(62,15).Value2 = 430000
$ws.Cells.Item(62,16).Value2 = 427500
$ws.Cells.Item(62,17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(62,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(62,19).Value2 = 1069
$ws.Cells.Item(62,20).Value2 = 400

# Row 63
$ws.Cells.Item(63,1).Value2 = $mercadoId
$ws.Cells.Item(63,2).Value2 = $mercado
$ws.Cells.Item(63,3).Value2 = $region
$ws.Cells.Item(63,4).Value2 = 44629
$ws.Cells.Item(63,5).Value2 = $codreg
$ws.Cells.Item(63,6).Value2 = $tipo
$ws.Cells.Item(63,7).Value2 = $productoId
$ws.Cells.Item(63,8).Value2 = $producto
$ws.Cells.Item(63,9).Value2 = $categoriaId
$ws.Cells.Item(63,10).Value2 = $categoria
$ws.Cells.Item(63,11).Value2 = "Phillips Cling"
$ws.Cells.Item(63,12).Value2 = "Primera"
$ws.Cells.Item(63,13).Value2 = 20
$ws.Cells.Item(63,14).Value2 = 385000
$ws.Cells.Item(63,15).Value2 = 390000
$ws.Cells.Item(63,16).Value2 = 387500
$ws.Cells.Item(63,17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(63,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(63,19).Value2 = 969
$ws.Cells.Item(63,20).Value2 = 400

# Row 64
$ws.Cells.Item(64,1).Value2 = $mercadoId
$ws.Cells.Item(64,2).Value2 = $mercado
$ws.Cells.Item(64,3).Value2 = $region
$ws.Cells.Item(64,4).Value2 = 44629
$ws.Cells.Item(64,5).Value2 = $codreg
$ws.Cells.Item(64,6).Value2 = $tipo
$ws.Cells.Item(64,7).Value2 = $productoId
$ws.Cells.Item(64,8).Value2 = $producto
$ws.Cells.Item(64,9).Value2 = $categoriaId
$ws.Cells.Item(64,10).Value2 = $categoria
$ws.Cells.Item(64,11).Value2 = "Phillips Cling"
$ws.Cells.Item(64,12).Value2 = "Segunda"
$ws.Cells.Item(64,13).Value2 = 16
$ws.Cells.Item(64,14).Value2 = 335000
$ws.Cells.Item(64,15).Value2 = 340000
$ws.Cells.Item(64,16).Value2 = 337500
$ws.Cells.Item(64,17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(64,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(64,19).Value2 = 844
$ws.Cells.Item(64,20).Value2 = 400

# Row 65
$ws.Cells.Item(65,1).Value2 = $mercadoId
$ws.Cells.Item(65,2).Value2 = $mercado
$ws.Cells.Item(65,3).Value2 = $region
$ws.Cells.Item(65,4).Value2 = 44629
$ws.Cells.Item(65,5).Value2 = $codreg
$ws.Cells.Item(65,6).Value2 = $tipo
$ws.Cells.Item(65,7).Value2 = $productoId
$ws.Cells.Item(65,8).Value2 = $producto
$ws.Cells.Item(65,9).Value2 = $categoriaId
$ws.Cells.Item(65,10).Value2 = $categoria
$ws.Cells.Item(65,11).Value2 = "September Snow"
$ws.Cells.Item(65,12).Value2 = "Especial"
$ws.Cells.Item(65,13).Value2 = 10
$ws.Cells.Item(65,14).Value2 = 400000
$ws.Cells.Item(65,15).Value2 = 410000
$ws.Cells.Item(65,16).Value2 = 405000
$ws.Cells.Item(65,17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(65,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(65,19).Value2 = 1012
$ws.Cells.Item(65,20).Value2 = 400

# Row 66
$ws.Cells.Item(66,1).Value2 = $mercadoId
$ws.Cells.Item(66,2).Value2 = $mercado
$ws.Cells.Item(66,3).Value2 = $region
$ws.Cells.Item(66,4).Value2 = 44629
$ws.Cells.Item(66,5).Value2 = $codreg
$ws.Cells.Item(66,6).Value2 = $tipo
$ws.Cells.Item(66,7).Value2 = $productoId
$ws.Cells.Item(66,8).Value2 = $producto
$ws.Cells.Item(66,9).Value2 = $categoriaId
$ws.Cells.Item(66,10).Value2 = $categoria
$ws.Cells.Item(66,11).Value2 = "September Snow"
$ws.Cells.Item(66,12).Value2 = "Primera"
$ws.Cells.Item(66,13).Value2 = 10
$ws.Cells.Item(66,14).Value2 = 340000
$ws.Cells.Item(66,15).Value2 = 350000
$ws.Cells.Item(66,16).Value2 = 345000
$ws.Cells.Item(66,17).Value2 = "$/bins (400 kilos)"
$ws.Cells.Item(66,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(66,19).Value2 = 862
$ws.Cells.Item(66,20).Value2 = 400
